$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.067.43"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.47"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.54"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.90"
$ws.Range("E6").Value = "  +2.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.110"
$ws.Range("E9").Value = "  +2.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.386"
$ws.Range("E10").Value = "  +6.66%  "

$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.62"
$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.114.44"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.837.05"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000148"
$ws.Range("E16").Value = "  +1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.641.93"
$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.82"
$ws.Range("E18").Value = "  +3.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("E19").Value = "  +3.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.02"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("E23").Value = "  -0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.42"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("E25").Value = "  +8.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.36"
$ws.Range("E26").Value = "  +7.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.70"
$ws.Range("E27").Value = "  +2.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "563.90"
$ws.Range("E28").Value = "  +4.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.18"
$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0859"
$ws.Range("E33").Value = "  +6.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.32"
$ws.Range("E35").Value = "  +2.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.32"
$ws.Range("E36").Value = "  -1.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.407"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.96"
$ws.Range("E39").Value = "  +4.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.25"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "165.47"
$ws.Range("E42").Value = "  -6.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.23"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.82"
$ws.Range("E44").Value = "  +1.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.09"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0572"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.629"
$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("E48").Value = "  +14.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0247"
$ws.Range("E49").Value = "  +2.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.93"
$ws.Range("E51").Value = "  +0.39%  "
